# Updated cryptos list on Sun May  7 02:51:44 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.014.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.99%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.904.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -4.21%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.21%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'324.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.06%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.00%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4597"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.46%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3807"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -3.04%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07710"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.80%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9730"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.11%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'21.90"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -4.36%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.931.76"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.47%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'6.927"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.69%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.649"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.52%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.07067"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.54%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +0.17%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'83.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -4.62%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000009487"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.67%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -4.13%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'1.005"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.19%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'28.977.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.16%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.298"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -4.46%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'10.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.94%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.100"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.37%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'158.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.64%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'19.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.51%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'5.585"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -3.96%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'117.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.68%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.847"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.68%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.09251"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.83%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.8571"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -4.32%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'5.081"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.91%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.238"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -7.26%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'2.975"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -6.88%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.05666"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.42%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.140"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.78%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.005"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.23%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.02034"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -3.09%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.5476"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -4.58%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'7.386"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -5.90%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1750"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -3.01%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'9.296"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -3.82%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.767"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.44%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.5153"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.91%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'11.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -7.04%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.077"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.24%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.06823"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.80%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.000002595"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -20.42%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'Quant"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'110.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -3.67%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'NEARProtocol"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'1.766"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.24%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.05%  "
$ws.Range("E51").Style = "Normal"
